{"js": "// Office.js (Word JavaScript API) edit script.\n// This is the body of `async (context) => { ... }`.\n//\n// Summary of the edit (see commit message \"Added many more features\"):\n//  1. Retitle the page (Heading1 + the repeated bold title near the bottom).\n//  2. Reword the first \"What we like\" bullet and add two brand-new bullets\n//     right after it (keeping \"Bonus round of free spins with variations\").\n//  3. Drop the \"Smooth gameplay and vibrant design\" / \"Suitable for casual\n//     players and experienced gamblers\" bullets from \"What we like\".\n//  4. Reword both \"What we don't like\" bullets.\n//  5. Reword the italic meta-description paragraph at the very end.\n\nconst body = context.document.body;\n\nconst OLD_TITLE = \"Play Day and Night Free Slot Game - Review 2021\";\nconst NEW_TITLE = \"Play Day and Night Free - Exciting Slot Game with Stunning Graphics\";\n\nconst OLD_WHAT_WE_LIKE_1 = \"Stunning graphics and game symbols\";\nconst NEW_WHAT_WE_LIKE_1 = \"High-quality graphics and game symbols\";\n\nconst NEW_BULLET_DEITIES = \"Well-designed and executed theme of deities\";\nconst NEW_BULLET_VISUALS = \"Stunning visuals and atmospheric music\";\n\nconst OLD_BULLET_SMOOTH = \"Smooth gameplay and vibrant design\";\nconst OLD_BULLET_SUITABLE = \"Suitable for casual players and experienced gamblers\";\n\nconst OLD_DONT_LIKE_1 = \"Low variance may not offer big winnings\";\nconst NEW_DONT_LIKE_1 = \"Bonus round cannot be reactivated\";\n\nconst OLD_DONT_LIKE_2 = \"Free spins cannot be reactivated\";\nconst NEW_DONT_LIKE_2 = \"High winnings are rare\";\n\nconst OLD_META_DESC =\n  \"Play Day and Night, an ancient Egyptian-themed slot game with stunning graphics and free spins. Read our review to learn more and play for free.\";\nconst NEW_META_DESC =\n  \"Read our review of Day and Night, an exciting slot game with stunning graphics. Play for free and enjoy the bonus round!\";\n\n// --- 1. Title (appears twice: the Heading1 and the bold \"recap\" line). -----\nconst titleHits = body.search(OLD_TITLE, { matchCase: true });\ntitleHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < titleHits.items.length; i++) {\n  titleHits.items[i].insertText(NEW_TITLE, \"Replace\");\n}\n\n// --- 2. First \"What we like\" bullet + two new bullets after it. -----------\nconst like1Hits = body.search(OLD_WHAT_WE_LIKE_1, { matchCase: true });\nlike1Hits.load(\"text\");\nawait context.sync();\nif (like1Hits.items.length > 0) {\n  const hit = like1Hits.items[0];\n  hit.insertText(NEW_WHAT_WE_LIKE_1, \"Replace\");\n  await context.sync();\n\n  const like1Para = hit.paragraphs.getFirst();\n  const deitiesPara = like1Para.insertParagraph(NEW_BULLET_DEITIES, \"After\");\n  deitiesPara.insertParagraph(NEW_BULLET_VISUALS, \"After\");\n}\n\n// --- 3. Remove the two obsolete \"What we like\" bullets. --------------------\nconst smoothHits = body.search(OLD_BULLET_SMOOTH, { matchCase: true });\nsmoothHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < smoothHits.items.length; i++) {\n  smoothHits.items[i].paragraphs.getFirst().delete();\n}\n\nconst suitableHits = body.search(OLD_BULLET_SUITABLE, { matchCase: true });\nsuitableHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < suitableHits.items.length; i++) {\n  suitableHits.items[i].paragraphs.getFirst().delete();\n}\n\n// --- 4. \"What we don't like\" bullets. --------------------------------------\nconst dontLike1Hits = body.search(OLD_DONT_LIKE_1, { matchCase: true });\ndontLike1Hits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < dontLike1Hits.items.length; i++) {\n  dontLike1Hits.items[i].insertText(NEW_DONT_LIKE_1, \"Replace\");\n}\n\nconst dontLike2Hits = body.search(OLD_DONT_LIKE_2, { matchCase: true });\ndontLike2Hits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < dontLike2Hits.items.length; i++) {\n  dontLike2Hits.items[i].insertText(NEW_DONT_LIKE_2, \"Replace\");\n}\n\n// --- 5. Italic meta-description paragraph. ---------------------------------\nconst metaHits = body.search(OLD_META_DESC, { matchCase: true });\nmetaHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < metaHits.items.length; i++) {\n  metaHits.items[i].insertText(NEW_META_DESC, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $d (ActiveDocument) are pre-seeded by the harness.\n#\n# Summary of the edit (see commit message \"Added many more features\"):\n#  1. Retitle the page (Heading1 + the repeated bold title near the bottom).\n#  2. Reword the first \"What we like\" bullet and add two brand-new bullets\n#     right after it (keeping \"Bonus round of free spins with variations\").\n#  3. Drop the \"Smooth gameplay and vibrant design\" / \"Suitable for casual\n#     players and experienced gamblers\" bullets from \"What we like\".\n#  4. Reword both \"What we don't like\" bullets.\n#  5. Reword the italic meta-description paragraph at the very end.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike, MatchAllWordForms,\n    # Forward, Wrap(wdFindContinue=1), Format, ReplaceWith, Replace(wdReplaceAll=2)\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n$OLD_TITLE = \"Play Day and Night Free Slot Game - Review 2021\"\n$NEW_TITLE = \"Play Day and Night Free - Exciting Slot Game with Stunning Graphics\"\n\n$OLD_WHAT_WE_LIKE_1 = \"Stunning graphics and game symbols\"\n$NEW_WHAT_WE_LIKE_1 = \"High-quality graphics and game symbols\"\n\n$NEW_BULLET_DEITIES = \"Well-designed and executed theme of deities\"\n$NEW_BULLET_VISUALS = \"Stunning visuals and atmospheric music\"\n\n$OLD_BULLET_SMOOTH = \"Smooth gameplay and vibrant design\"\n$OLD_BULLET_SUITABLE = \"Suitable for casual players and experienced gamblers\"\n\n$OLD_DONT_LIKE_1 = \"Low variance may not offer big winnings\"\n$NEW_DONT_LIKE_1 = \"Bonus round cannot be reactivated\"\n\n$OLD_DONT_LIKE_2 = \"Free spins cannot be reactivated\"\n$NEW_DONT_LIKE_2 = \"High winnings are rare\"\n\n$OLD_META_DESC = \"Play Day and Night, an ancient Egyptian-themed slot game with stunning graphics and free spins. Read our review to learn more and play for free.\"\n$NEW_META_DESC = \"Read our review of Day and Night, an exciting slot game with stunning graphics. Play for free and enjoy the bonus round!\"\n\n# --- 2. First \"What we like\" bullet + two new bullets right after it. ------\n# Find the paragraph that currently holds the bullet, insert the two new\n# list-bulleted paragraphs after it (inheriting its pPr/style), then reword\n# the original bullet's text in place.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $OLD_WHAT_WE_LIKE_1) {\n        $p.Range.InsertParagraphAfter()\n        $deitiesPara = $p.Next()\n        $deitiesPara.Range.Text = $NEW_BULLET_DEITIES\n\n        $deitiesPara.Range.InsertParagraphAfter()\n        $visualsPara = $deitiesPara.Next()\n        $visualsPara.Range.Text = $NEW_BULLET_VISUALS\n        break\n    }\n}\n\n# --- 3. Remove the two obsolete \"What we like\" bullets. --------------------\nforeach ($targetText in @($OLD_BULLET_SMOOTH, $OLD_BULLET_SUITABLE)) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText) {\n            $p.Range.Delete() | Out-Null\n            break\n        }\n    }\n}\n\n# --- 1, 2 (reword), 4, 5: plain text replacements (global, exact match). ---\nReplace-AllText $OLD_TITLE $NEW_TITLE\nReplace-AllText $OLD_WHAT_WE_LIKE_1 $NEW_WHAT_WE_LIKE_1\nReplace-AllText $OLD_DONT_LIKE_1 $NEW_DONT_LIKE_1\nReplace-AllText $OLD_DONT_LIKE_2 $NEW_DONT_LIKE_2\nReplace-AllText $OLD_META_DESC $NEW_META_DESC\n\n\"done\"\n"}
